# Arreglos generales en scripts de generación de siniestros
#
# Updates the two "i-preproducciongestion" rows (claims #1 and #2) to point at
# the "ssurgwsoadev4-oci" environment/URL used by the rest of the sheet, gives
# them fresh policy numbers and claim dates, and tweaks the "Descripcion" on
# the first row. Also removes the now-stale hyperlink object that covered
# C2:C8 (since C2/C3 are no longer links to the old ClaimCenter URL) and
# updates the saved window view (scroll position / selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (claim #2) -------------------------------------------------
# Written before row 2 purely so the shared-string table grows in the same
# order as the reference edit (cosmetic, but keeps the xlsx byte-for-byte
# closer to the target).
$ws.Range("F3").Value = "04104015431"
$ws.Range("H3").Value = "'21/03/2022"
$ws.Range("B3").Value = "ssurgwsoadev4-oci.opc.oracleoutsourcing.com"

# --- Row 2 (claim #1) ---------------------------------------------------
$ws.Range("F2").Value = "04104015535"
$ws.Range("H2").Value = "'28/03/2023"
$ws.Range("B2").Value = "ssurgwsoadev4-oci.opc.oracleoutsourcing.com"
$ws.Range("U2").Value = "Cristales"

# --- Remove the stale C2:C8 hyperlink and retarget C2/C3 ----------------
$staleLink = $null
for ($i = 1; $i -le $ws.Hyperlinks.Count; $i++) {
    $h = $ws.Hyperlinks.Item($i)
    if ($h.Range.Address() -eq '$C$2:$C$8') {
        $staleLink = $h
    }
}
if ($staleLink -ne $null) {
    $staleLink.Delete()
}

$ws.Range("C2:C3").Style = "Normal"
$ws.Range("C2").Value = "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/cc/ClaimCenter.do"
$ws.Range("C3").Value = "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/cc/ClaimCenter.do"

# --- Saved window view: scroll back to A1, select H3 --------------------
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H3").Select() | Out-Null
